$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.492.55"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "2.428.98"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "504.85"
$ws.Range("E5").Value = "  -3.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.51"
$ws.Range("E6").Value = "  -3.75%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -2.26%  "
$ws.Range("D9").Value = "2.440.27"
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0946"
$ws.Range("E11").Value = "  -5.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.15"
$ws.Range("E12").Value = "  -4.18%  "
$ws.Range("E13").Value = "  -5.42%  "
$ws.Range("D14").Value = "2.861.12"
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("D15").Value = "57.425.04"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.63"
$ws.Range("E16").Value = "  -2.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000132"
$ws.Range("E17").Value = "  -3.79%  "
$ws.Range("D18").Value = "2.438.57"
$ws.Range("E18").Value = "  -2.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.39"
$ws.Range("E19").Value = "  -4.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "315.94"
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.08"
$ws.Range("E21").Value = "  -2.83%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.72"
$ws.Range("E23").Value = "  -1.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.79"
$ws.Range("E24").Value = "  -2.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.402"
$ws.Range("E25").Value = "  -2.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.15"
$ws.Range("E28").Value = "  -3.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "168.41"
$ws.Range("E29").Value = "  +0.60%  "
$ws.Range("D30").Value = "0.0₃0717"
$ws.Range("E30").Value = "  -4.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.17"
$ws.Range("E31").Value = "  -3.17%  "
$ws.Range("E32").Value = "  -3.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.12"
$ws.Range("E33").Value = "  -7.50%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.62"
$ws.Range("E36").Value = "  -3.11%  "
$ws.Range("E37").Value = "  -6.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.85"
$ws.Range("E38").Value = "  -4.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.41"
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.44"
$ws.Range("E40").Value = "  -3.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.747"
$ws.Range("E41").Value = "  -6.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "268.41"
$ws.Range("E42").Value = "  -3.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.34"
$ws.Range("E43").Value = "  -4.31%  "
$ws.Range("E44").Value = "  -5.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.577"
$ws.Range("E45").Value = "  -3.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0907"
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "118.25"
$ws.Range("E47").Value = "  -6.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0482"
$ws.Range("E48").Value = "  -2.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.99"
$ws.Range("E49").Value = "  -5.45%  "
$ws.Range("E50").Value = "  -4.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.41"
$ws.Range("E51").Value = "  -4.98%  "
